$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D-column cells whose new value would otherwise be
# auto-converted to a number by Excel, so the stored text matches the source
# exactly (e.g. "1.00" must stay "1.00", not become 1).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update cell values per row
$ws.Range("D2").Value = "96.743.55"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "3.701.43"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").Value = "2.44"
$ws.Range("E4").Value = "  +29.99%  "

$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "229.49"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("D7").Value = "653.09"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").Value = "0.439"
$ws.Range("E8").Value = "  +2.59%  "

$ws.Range("D9").Value = "1.16"
$ws.Range("E9").Value = "  +10.02%  "

$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "3.699.95"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "47.83"
$ws.Range("E12").Value = "  +8.10%  "

$ws.Range("D13").Value = "0.211"
$ws.Range("E13").Value = "  +3.05%  "

$ws.Range("D14").Value = "0.0000301"
$ws.Range("E14").Value = "  -2.55%  "

$ws.Range("D15").Value = "6.64"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").Value = "4.394.92"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").Value = "96.589.32"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("E18").Value = "  +1.03%  "

$ws.Range("D19").Value = "3.700.98"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").Value = "19.69"
$ws.Range("E20").Value = "  +5.40%  "

$ws.Range("D21").Value = "12.96"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("E22").Value = "  +9.03%  "

$ws.Range("D23").Value = "534.99"
$ws.Range("E23").Value = "  +3.31%  "

$ws.Range("D24").Value = "3.33"
$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("D25").Value = "0.253"
$ws.Range("E25").Value = "  +43.51%  "

$ws.Range("D26").Value = "120.94"
$ws.Range("E26").Value = "  +20.22%  "

$ws.Range("D27").Value = "0.0000210"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").Value = "3.902.56"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").Value = "13.02"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").Value = "13.37"
$ws.Range("E31").Value = "  +10.62%  "

$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D34").Value = "0.188"
$ws.Range("E34").Value = "  +1.97%  "

$ws.Range("D35").Value = "33.40"
$ws.Range("E35").Value = "  +4.10%  "

$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("E37").Value = "  +4.73%  "

$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "609.81"
$ws.Range("E39").Value = "  -6.92%  "

$ws.Range("D41").Value = "8.46"
$ws.Range("E41").Value = "  -3.62%  "

$ws.Range("D42").Value = "7.19"
$ws.Range("E42").Value = "  +3.55%  "

$ws.Range("D43").Value = "0.510"
$ws.Range("E43").Value = "  +18.83%  "

$ws.Range("D44").Value = "0.0509"
$ws.Range("E44").Value = "  +13.31%  "

$ws.Range("D45").Value = "0.163"
$ws.Range("E45").Value = "  +2.77%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "40.73"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -3.67%  "

$ws.Range("D48").Value = "0.976"
$ws.Range("E48").Value = "  +2.49%  "

$ws.Range("D49").Value = "9.04"
$ws.Range("E49").Value = "  +7.13%  "

$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").Value = "23.55"
$ws.Range("E51").Value = "  -0.04%  "
